# Apply updated crafting profit values per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2016494.2
$ws.Range("J38").Value = 609.875
$ws.Range("L38").Value = 1829.625
$ws.Range("N38").Value = -2573.625
$ws.Range("H96").Value = 803.2
$ws.Range("I96").Value = 607.7
$ws.Range("J96").Value = 1194.2
$ws.Range("K96").Value = 1823.1
$ws.Range("L96").Value = 3582.6
$ws.Range("M96").Value = -450.1000000000001
$ws.Range("N96").Value = -6328.6
$ws.Range("H116").Value = 2276.5
$ws.Range("I116").Value = 1601.6666
$ws.Range("J116").Value = 2565.7144
$ws.Range("K116").Value = 1601.6666
$ws.Range("L116").Value = 2565.7144
$ws.Range("M116").Value = 1840.3334
$ws.Range("N116").Value = -9449.714400000001
$ws.Range("H131").Value = 4435.6
$ws.Range("I131").Value = 2802.5
$ws.Range("J131").Value = 4686.846
$ws.Range("K131").Value = 8407.5
$ws.Range("L131").Value = 14060.538
$ws.Range("M131").Value = -3367.5
$ws.Range("N131").Value = -24140.538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H45").Value = 1752.48
$ws.Range("I45").Value = 1681.2142
$ws.Range("J45").Value = 1843.1818
$ws.Range("K45").Value = 1681.2142
$ws.Range("L45").Value = 1843.1818
$ws.Range("M45").Value = -1304.2142
$ws.Range("N45").Value = -2597.1818
$ws.Range("H80").Value = 13905
$ws.Range("J80").Value = 14556
$ws.Range("L80").Value = 14556
$ws.Range("N80").Value = -16552
$ws.Range("H83").Value = 13905
$ws.Range("J83").Value = 14556
$ws.Range("L83").Value = 43668
$ws.Range("N83").Value = -53652
$ws.Range("H122").Value = 1973.4231
$ws.Range("I122").Value = 1938.7222
$ws.Range("K122").Value = 5816.1666
$ws.Range("M122").Value = -3366.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 933.3333
$ws.Range("I6").Value = 500
$ws.Range("K6").Value = 500
$ws.Range("M6").Value = -387
$ws.Range("H7").Value = 201.41176
$ws.Range("I7").Value = 62
$ws.Range("J7").Value = 457
$ws.Range("K7").Value = 62
$ws.Range("L7").Value = 457
$ws.Range("M7").Value = 51
$ws.Range("N7").Value = -683
$ws.Range("H68").Value = 12903.546
$ws.Range("J68").Value = 12903.546
$ws.Range("L68").Value = 12903.546
$ws.Range("N68").Value = -14401.546
$ws.Range("H71").Value = 12903.546
$ws.Range("J71").Value = 12903.546
$ws.Range("L71").Value = 38710.638
$ws.Range("N71").Value = -46198.638
$ws.Range("H74").Value = 19952.75
$ws.Range("J74").Value = 19952.75
$ws.Range("L74").Value = 19952.75
$ws.Range("N74").Value = -21700.75
$ws.Range("H77").Value = 19952.75
$ws.Range("J77").Value = 19952.75
$ws.Range("L77").Value = 59858.25
$ws.Range("N77").Value = -68594.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1789.8
$ws.Range("J34").Value = 2199.75
$ws.Range("L34").Value = 6599.25
$ws.Range("N34").Value = -6767.25
$ws.Range("H39").Value = 22200
$ws.Range("J39").Value = 22200
$ws.Range("L39").Value = 66600
$ws.Range("N39").Value = -67188
$ws.Range("H55").Value = 8413.588
$ws.Range("J55").Value = 8901.9375
$ws.Range("L55").Value = 26705.8125
$ws.Range("N55").Value = -27059.8125
$ws.Range("H68").Value = 1806.8986
$ws.Range("I68").Value = 1147.4688
$ws.Range("J68").Value = 2377.2163
$ws.Range("K68").Value = 3442.4064
$ws.Range("L68").Value = 7131.6489
$ws.Range("M68").Value = -2631.4064
$ws.Range("N68").Value = -8753.6489
$ws.Range("H71").Value = 1806.8986
$ws.Range("I71").Value = 1147.4688
$ws.Range("J71").Value = 2377.2163
$ws.Range("K71").Value = 10327.2192
$ws.Range("L71").Value = 21394.9467
$ws.Range("M71").Value = -6271.219200000001
$ws.Range("N71").Value = -29506.9467
$ws.Range("H132").Value = 3113.611
$ws.Range("I132").Value = 4170
$ws.Range("J132").Value = 2902.3333
$ws.Range("K132").Value = 37530
$ws.Range("L132").Value = 26120.9997
$ws.Range("M132").Value = -35000
$ws.Range("N132").Value = -31180.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 12699
$ws.Range("J15").Value = 12699
$ws.Range("L15").Value = 12699
$ws.Range("N15").Value = -13275
$ws.Range("H43").Value = 4683.7
$ws.Range("J43").Value = 5862.4287
$ws.Range("L43").Value = 5862.4287
$ws.Range("N43").Value = -6164.4287
$ws.Range("H58").Value = 17524.75
$ws.Range("J58").Value = 22699.666
$ws.Range("L58").Value = 22699.666
$ws.Range("N58").Value = -23253.666
$ws.Range("H81").Value = 12699
$ws.Range("J81").Value = 12699
$ws.Range("L81").Value = 12699
$ws.Range("N81").Value = -14695
$ws.Range("H84").Value = 12699
$ws.Range("J84").Value = 12699
$ws.Range("L84").Value = 38097
$ws.Range("N84").Value = -48081
$ws.Range("H107").Value = 459396.5
$ws.Range("I107").Value = 198.58333
$ws.Range("J107").Value = 1010434
$ws.Range("K107").Value = 198.58333
$ws.Range("L107").Value = 1010434
$ws.Range("M107").Value = 1721.41667
$ws.Range("N107").Value = -1014274
$ws.Range("H122").Value = 3708.9375
$ws.Range("I122").Value = 4019.1
$ws.Range("J122").Value = 3192
$ws.Range("K122").Value = 12057.3
$ws.Range("L122").Value = 9576
$ws.Range("M122").Value = -9607.299999999999
$ws.Range("N122").Value = -14476

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 563394.4
$ws.Range("J46").Value = 1266999.9
$ws.Range("L46").Value = 1266999.9
$ws.Range("N46").Value = -1267375.9
$ws.Range("H55").Value = 516969.22
$ws.Range("I55").Value = 874303.4399999999
$ws.Range("J55").Value = 819.7778
$ws.Range("K55").Value = 874303.4399999999
$ws.Range("L55").Value = 819.7778
$ws.Range("M55").Value = -874130.4399999999
$ws.Range("N55").Value = -1165.7778
$ws.Range("H56").Value = 20000
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 20000
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 20000
$ws.Range("N56").Value = -21382
$ws.Range("M56").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 17949
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 17949
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 17949
$ws.Range("N61").Value = -18533
$ws.Range("M61").ClearContents()
$ws.Range("H96").Value = 66668604
$ws.Range("I96").Value = 125001970
$ws.Range("J96").Value = 1902.4286
$ws.Range("K96").Value = 125001970
$ws.Range("L96").Value = 1902.4286
$ws.Range("M96").Value = -125000597
$ws.Range("N96").Value = -4648.4286
